$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header style (from H1, which carries the bold/centered/
# bordered header style) onto the two new header cells so they match the
# look of the rest of row 1.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

# New headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for columns I (I0) and J (IF), rows 2-30
$values = @(
    @(1, 5),
    @(1, 5),
    @(1, 5),
    @(1, 7),
    @(1, 6),
    @(1, 6),
    @(1, 6),
    @(1, 7),
    @(1, 6),
    @(1, 6),
    @(1, 5),
    @(1, 5),
    @(6, 7),
    @(1, 5),
    @(1, 4),
    @(6, 7),
    @(9, 9),
    @(7, 9),
    @(5, 7),
    @(1, 4),
    @(8, 9),
    @(3, 6),
    @(9, 9),
    @(7, 8),
    @(9, 9),
    @(7, 7),
    @(7, 7),
    @(7, 7),
    @(6, 6)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $pair = $values[$i]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
